$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly price log for "Betarraga" at "Terminal La Palmera de
# La Serena". Every week two rows are prepended (a "Primera" and a "Segunda"
# grade row) and the oldest week's two rows fall off... except here a new
# week was inserted at the top (rows 34/35 got a newer date) which pushes
# every following pair of rows (36..157) down by one pair (two rows), and
# the two rows that used to be at the very bottom (156/157) reappear as two
# brand new rows (158/159) at the end.
#
# Columns that vary per-row: D (Fecha), J/K/L/M (prices), P. Columns
# A,B,C,E,F,G,H,I,N,O,Q,R only depend on whether the row is "Primera" or
# "Segunda" and are already correct/constant, so they don't need touching.

$valueCols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P
$allCols   = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)  # A..R

# 1) Capture the two rows that are about to be "pushed off the end" (156 and
#    157) BEFORE anything is overwritten, so we can re-create them as the new
#    rows 158/159.
$saved156 = @{}
$saved157 = @{}
foreach ($c in $allCols) {
    $saved156[$c] = $ws.Cells.Item(156, $c).Value()
    $saved157[$c] = $ws.Cells.Item(157, $c).Value()
}

# 2) Shift rows 36..157 down by one pair: row r takes the old values that
#    used to live in row (r - 2). Walk from the bottom up so we never
#    clobber a row before it has been read as a source.
for ($r = 157; $r -ge 36; $r--) {
    $src = $r - 2
    foreach ($c in $valueCols) {
        $v = $ws.Cells.Item($src, $c).Value()
        if ($c -eq 4) {
            # Keep the date column's existing custom date format instead of
            # letting a fresh auto-detected date style get created.
            $ws.Cells.Item($r, $c).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        }
        $ws.Cells.Item($r, $c).Value = $v
    }
}

# 3) Re-create the old rows 156/157 as new rows 158/159 (full row copy).
foreach ($c in $allCols) {
    if ($c -eq 4) {
        $ws.Cells.Item(158, $c).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item(158, $c).Value = $saved156[$c]
}
foreach ($c in $allCols) {
    if ($c -eq 4) {
        $ws.Cells.Item(159, $c).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item(159, $c).Value = $saved157[$c]
}

# 4) The newest week (rows 34/35) gets its own new date/price instead of
#    being shifted from somewhere else.
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 4).Value = 44453
$ws.Cells.Item(34, 10).Value = 3400

$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 4).Value = 44453
